# Fix a deckview show bug on card level up.
# The "TowerLevel" column (D) on the Exp sheet was giving every card level
# from 18..99 (rows 21..102) a flat TowerLevel of 3, which under-reported the
# tower level needed for higher level cards and caused the deck view to show
# the wrong state when a card leveled up. Re-stage the TowerLevel curve so it
# increases every ~7-9 rows as the player level goes up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exp")

$ws.Range("D21:D27").Value = 4
$ws.Range("D28:D36").Value = 5
$ws.Range("D37:D48").Value = 6
$ws.Range("D49:D62").Value = 7
$ws.Range("D63:D85").Value = 8
$ws.Range("D86:D101").Value = 9
$ws.Range("D102").Value = 10

# Leave the cursor where the author last looked while checking the fix.
$ws.Range("D101").Select()
